$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")
$ws.Activate()

# Candidate power plants: switch "realistic capacity" flags from FALSE to TRUE
$ws.Range("B18").Value = $true
$ws.Range("B19").Value = $true

# Update descriptive formulas (B18/B19 change also recalculates C20 automatically)
$ws.Range("C25").Formula = '=IF(B24=FALSE,"- > NOT ACTIVE, prices are not being fixed, to do so change previous like to TRUE","fixed prices for investment")'
$ws.Range("C27").Formula = '=IF(B27=FALSE,"profiles for the ACTUAL year but not for the future year (investors see only one year)"," ")'

# Restore the view / selection state as it was left by the author
$ws.Range("C27").Select()
